{"js": "// Adds two new paragraphs at the very end of the document body (after the\n// existing trailing empty paragraph, right before the section break):\n//   \"Paper 1: \" + <title of paper 1>\n//   \"Paper 2: \" + <title of paper 2>\n// Each paragraph is made of two runs (label run + title run), matching the\n// target OOXML exactly.\n\n// Minimal XML-escaping helper so the label/title text is safe to splice into\n// a raw OOXML (flat-OPC) fragment.\nfunction xmlEscape(s) {\n  return String(s)\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\n// Builds a flat-OPC package wrapping a single <w:p> with two runs: a\n// \"label\" run (kept with xml:space=\"preserve\" since it ends in a space)\n// and a \"title\" run, exactly mirroring the target diff's XML shape.\nfunction singleParagraphFlatOpc(label, title) {\n  const labelXml = xmlEscape(label);\n  const titleXml = xmlEscape(title);\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    \"<w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">' + labelXml + \"</w:t></w:r>\" +\n    \"<w:r><w:t>\" + titleXml + \"</w:t></w:r>\" +\n    \"</w:p>\" +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// Appends a new two-run paragraph at the end of the body, leaving any\n// existing trailing (e.g. empty) paragraph untouched before it.\nfunction appendTwoRunParagraph(body, label, title) {\n  // Insert a fresh empty paragraph at the very end, then stamp its OOXML\n  // (via \"Replace\") so it ends up with exactly two <w:r> runs instead of\n  // Word's usual same-format run coalescing that plain insertText() would\n  // trigger.\n  const newPara = body.insertParagraph(\"\", \"End\");\n  newPara.insertOoxml(singleParagraphFlatOpc(label, title), \"Replace\");\n}\n\nconst body = context.document.body;\n\nappendTwoRunParagraph(\n  body,\n  \"Paper 1: \",\n  \"Diversity of fishing gears and crafts used for harvesting the Asian seabass, Lates calcarifer along the Bay of Bengal, Bangladesh coast\"\n);\n\nappendTwoRunParagraph(\n  body,\n  \"Paper 2: \",\n  \"Fishing Gears and Practices in the Bukbhora Oxbow Lake: Implications for Biodiversity Conservation in South-west Bangladesh\"\n);\n\nawait context.sync();\n", "ps1": "# Adds two new paragraphs at the very end of the document body (after the\n# existing trailing empty paragraph, right before the section break):\n#   \"Paper 1: \" + <title of paper 1>\n#   \"Paper 2: \" + <title of paper 2>\n# Each paragraph is made of two runs (a \"label\" run and a \"title\" run),\n# matching the target OOXML exactly (two separate <w:r> elements rather\n# than one run Word would otherwise coalesce them into).\n\nfunction Escape-Xml([string]$s) {\n    if ($null -eq $s) { return \"\" }\n    return $s.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\").Replace('\"', \"&quot;\").Replace(\"'\", \"&apos;\")\n}\n\nfunction New-FlatOpcParagraph([string]$label, [string]$title) {\n    $labelXml = Escape-Xml $label\n    $titleXml = Escape-Xml $title\n    return @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">$labelXml</w:t></w:r>\n            <w:r><w:t>$titleXml</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n}\n\n# Appends a new two-run paragraph at the end of the document, leaving any\n# existing trailing (e.g. empty) paragraph untouched before it.\nfunction Add-TwoRunParagraph($d, [string]$label, [string]$title) {\n    $lastPara = $d.Paragraphs.Last\n    $lastPara.Range.InsertParagraphAfter()\n\n    # The freshly inserted (still empty) paragraph is a single-character\n    # range holding just the paragraph mark; collapse to its end (right\n    # before that mark) so InsertXML's own <w:p>...</w:p> mark replaces it,\n    # landing our two runs into a brand new paragraph.\n    $newPara = $d.Paragraphs.Last\n    $insertionRange = $newPara.Range\n    $insertionRange.Collapse(0)  # wdCollapseEnd\n    $insertionRange.InsertXML((New-FlatOpcParagraph $label $title))\n\n    # InsertXML's incoming paragraph mark pushes a fresh empty paragraph\n    # mark after it (Word always needs a terminating mark), so the story\n    # now has one extra blank paragraph at the tail. Fold it back into the\n    # paragraph we just created by deleting the paragraph mark between them.\n    $strayIndex = $d.Paragraphs.Count\n    $stray = $d.Paragraphs($strayIndex)\n    $content = $d.Paragraphs($strayIndex - 1)\n    if ($stray.Range.Text.Trim() -eq \"\") {\n        $mark = $d.Range($content.Range.End - 1, $stray.Range.End - 1)\n        $mark.Delete()\n    }\n}\n\n$d = $word.ActiveDocument\n\nAdd-TwoRunParagraph $d \"Paper 1: \" \"Diversity of fishing gears and crafts used for harvesting the Asian seabass, Lates calcarifer along the Bay of Bengal, Bangladesh coast\"\nAdd-TwoRunParagraph $d \"Paper 2: \" \"Fishing Gears and Practices in the Bukbhora Oxbow Lake: Implications for Biodiversity Conservation in South-west Bangladesh\"\n"}
